# Apply theoretical-framework notes update to Notas.docx
# Splits single-run paragraphs into multiple runs wrapped with w:proofErr
# spell-check markers around English/foreign terms and proper nouns, and
# appends the new "Scraps" section with the dominance/iterated-reasoning text.
$d = $word.ActiveDocument

# --- Paragraph 1: outlier / Boxplot / plot ---
$d.Paragraphs(1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Hay que demostrar que Participante A de la sesión 3 es un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>outlier</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. ¿</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Boxplot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">? Graficar todas las tiradas en un solo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>plot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (?).</w:t></w:r></w:p>')

# --- Paragraph 2: Nagel / Crowford / Camerer ---
$d.Paragraphs(2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">La introducción debe ir del pensamiento iterado a Keynes, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Nagel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, describir el juego, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Crowford</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Camerer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')

# --- Paragraph 3: Lahav ---
$d.Paragraphs(3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">En </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Lahav</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, contrario a lo que dice el artículo, se encontró Creencia &gt; Elección &gt; Elección*p</w:t></w:r></w:p>')

# --- Paragraph 4: Slonim ---
$d.Paragraphs(4).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Slonim</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> usó la mediana en el juego, yo usé la media. Buscar e incluir referencia sobre cómo afecta el juego usar una u otra.</w:t></w:r></w:p>')

# --- Paragraph 5: Lahav (normalized measure) ---
$d.Paragraphs(5).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Se desestima utilizar la medida normalizada de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Lahav</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> porque las diferencias se ponderan de forma arbitraria en cada periodo.</w:t></w:r></w:p>')

# --- Paragraph 6: Lahav (BO reporting) ---
$d.Paragraphs(6).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Vamos a dejar de utilizar BA y BA*p. Reportar BO solo al inicio igual que </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Lahav</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, aunque la medida importante es BO*p.</w:t></w:r></w:p>')

# --- Paragraph 7 onward: trailing space tweak on the last bullet, drop the
#     now-orphaned " " run + _GoBack bookmark from it, and append the new
#     "HAY QUE INVESTIGAR..." bullet, the "Scraps:" bullet and the six
#     indented paragraphs of theory notes (bookmark now sits on the final one).
$d.Paragraphs(7).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Graficar frecuencias relativas (sobre 100) de elecciones en el primer periodo. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>HAY QUE INVESTIGAR SI MEJORAN LAS CREENCIAS</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Scraps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>Uno de los conceptos básicos en teoría de juegos es la dominancia. Una estrategia A domina a una estrategia B cuando todos los pagos (i. e. consecuencias) asociados con la estrategia A son mejores que los asociados con la estrategia B, con independencia de las estrategias que puedan utilizar otros jugadores.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>En un juego donde existen estrategias que son dominadas por otras, las estrategias dominadas deben descartarse, lo que puede cambiar la estructura del juego ya que otras estrategias se convierten en estrategias dominadas. Estas también deben ser descartadas.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>Un juego tiene una solución por dominancia si pueden eliminarse estrategias dominadas repetidamente hasta que únicamente quede una estrategia dominante para cada jugador (i. e. el equilibrio de Nash).</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>Este proceso de eliminación de estrategias dominadas se conoce como razonamiento iterado.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>Resulta inverosímil pensar que las personas pueden repetir el razonamiento iterado tantas veces como sea necesario para llegar al equilibrio (en teoría, hasta una cantidad infinita de veces), dado que existen límites en las capacidades cognitivas de las personas. Hay evidencia de que las personas no empiezan jugando en equilibrio en juegos con muchas estrategias que tienen solución por dominancia (Stahl &amp; Wilson, 1995).</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:t>Si un jugador A basa su elección en las creencias que tiene sobre la elección que va a realizar un jugador B, el jugador A está utilizando un paso de razonamiento iterado. Adicionalmente, el jugador A puede preguntarse si el jugador B también está utilizando el mismo tipo de razonamiento, e incorporar esta creencia en su toma de decisiones. En este caso, el jugador A estaría usando dos pasos de razonamiento.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>El jugador A puede repetir este razonamiento una y otra vez (de ahí que se llame razonamiento iterado), incorporando en su elección las creencias que tiene sobre cuántas veces va a realizar este razonamiento el jugador B.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

